# Auto-generated edit script: updates market-price derived columns (H-N)
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect
# a refreshed data pull from the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- ALC sheet: 61 cell update(s) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 297.35715
$ws.Cells.Item(33, 9).Value = 273.92307
$ws.Cells.Item(33, 10).Value = 602
$ws.Cells.Item(33, 11).Value = 273.92307
$ws.Cells.Item(33, 12).Value = 602
$ws.Cells.Item(33, 13).Value = -44.92307
$ws.Cells.Item(33, 14).Value = -1060
$ws.Cells.Item(106, 8).Value = 3000
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 13).Value = $null
$ws.Cells.Item(108, 8).Value = 39560
$ws.Cells.Item(108, 10).Value = 39560
$ws.Cells.Item(108, 12).Value = 39560
$ws.Cells.Item(108, 14).Value = -47240
$ws.Cells.Item(113, 8).Value = 4940.1665
$ws.Cells.Item(113, 9).Value = 2409.75
$ws.Cells.Item(113, 10).Value = 10001
$ws.Cells.Item(113, 11).Value = 2409.75
$ws.Cells.Item(113, 12).Value = 10001
$ws.Cells.Item(113, 13).Value = 844.25
$ws.Cells.Item(113, 14).Value = -16509
$ws.Cells.Item(115, 8).Value = 1027.75
$ws.Cells.Item(115, 9).Value = 1027.75
$ws.Cells.Item(115, 11).Value = 3083.25
$ws.Cells.Item(115, 13).Value = -1516.25
$ws.Cells.Item(124, 8).Value = 41890
$ws.Cells.Item(124, 10).Value = 41890
$ws.Cells.Item(124, 12).Value = 41890
$ws.Cells.Item(124, 14).Value = -51710
$ws.Cells.Item(126, 8).Value = 41890
$ws.Cells.Item(126, 10).Value = 41890
$ws.Cells.Item(126, 12).Value = 41890
$ws.Cells.Item(126, 14).Value = -51770
$ws.Cells.Item(127, 8).Value = 2390.1
$ws.Cells.Item(127, 9).Value = 480.66666
$ws.Cells.Item(127, 11).Value = 1441.99998
$ws.Cells.Item(127, 13).Value = 3518.00002
$ws.Cells.Item(128, 8).Value = 41890
$ws.Cells.Item(128, 10).Value = 41890
$ws.Cells.Item(128, 12).Value = 41890
$ws.Cells.Item(128, 14).Value = -51850
$ws.Cells.Item(129, 8).Value = 823.09
$ws.Cells.Item(129, 10).Value = 867.61957
$ws.Cells.Item(129, 12).Value = 2602.85871
$ws.Cells.Item(129, 14).Value = -12602.85871
$ws.Cells.Item(130, 8).Value = 41853.332
$ws.Cells.Item(130, 10).Value = 41853.332
$ws.Cells.Item(130, 12).Value = 41853.332
$ws.Cells.Item(130, 14).Value = -51893.332
$ws.Cells.Item(133, 8).Value = 45000
$ws.Cells.Item(133, 10).Value = 45000
$ws.Cells.Item(133, 12).Value = 45000
$ws.Cells.Item(133, 14).Value = -55120
$ws.Cells.Item(141, 8).Value = 183962
$ws.Cells.Item(141, 9).Value = 223943
$ws.Cells.Item(141, 10).Value = 4047.5
$ws.Cells.Item(141, 11).Value = 671829
$ws.Cells.Item(141, 12).Value = 12142.5
$ws.Cells.Item(141, 13).Value = -666649
$ws.Cells.Item(141, 14).Value = -22502.5

# --- ARM sheet: 46 cell update(s) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5653.983
$ws.Cells.Item(32, 9).Value = 4118.3096
$ws.Cells.Item(32, 10).Value = 9448
$ws.Cells.Item(32, 11).Value = 4118.3096
$ws.Cells.Item(32, 12).Value = 9448
$ws.Cells.Item(32, 13).Value = -3831.3096
$ws.Cells.Item(32, 14).Value = -10022
$ws.Cells.Item(63, 8).Value = 8151351.5
$ws.Cells.Item(63, 9).Value = 11545248
$ws.Cells.Item(63, 10).Value = 6000
$ws.Cells.Item(63, 11).Value = 11545248
$ws.Cells.Item(63, 12).Value = 6000
$ws.Cells.Item(63, 13).Value = -11544562
$ws.Cells.Item(63, 14).Value = -7372
$ws.Cells.Item(66, 8).Value = 8151351.5
$ws.Cells.Item(66, 9).Value = 11545248
$ws.Cells.Item(66, 10).Value = 6000
$ws.Cells.Item(66, 11).Value = 57726240
$ws.Cells.Item(66, 12).Value = 30000
$ws.Cells.Item(66, 13).Value = -57722808
$ws.Cells.Item(66, 14).Value = -36864
$ws.Cells.Item(74, 8).Value = 2240.8462
$ws.Cells.Item(74, 9).Value = 803.1667
$ws.Cells.Item(74, 10).Value = 3473.1428
$ws.Cells.Item(74, 11).Value = 803.1667
$ws.Cells.Item(74, 12).Value = 3473.1428
$ws.Cells.Item(74, 13).Value = 70.83330000000001
$ws.Cells.Item(74, 14).Value = -5221.1428
$ws.Cells.Item(77, 8).Value = 2240.8462
$ws.Cells.Item(77, 9).Value = 803.1667
$ws.Cells.Item(77, 10).Value = 3473.1428
$ws.Cells.Item(77, 11).Value = 4015.8335
$ws.Cells.Item(77, 12).Value = 17365.714
$ws.Cells.Item(77, 13).Value = 352.1665000000003
$ws.Cells.Item(77, 14).Value = -26101.714
$ws.Cells.Item(132, 8).Value = 3395.5833
$ws.Cells.Item(132, 9).Value = 1878.25
$ws.Cells.Item(132, 10).Value = 4154.25
$ws.Cells.Item(132, 11).Value = 5634.75
$ws.Cells.Item(132, 12).Value = 12462.75
$ws.Cells.Item(132, 13).Value = -3104.75
$ws.Cells.Item(132, 14).Value = -17522.75
$ws.Cells.Item(137, 8).Value = 38568
$ws.Cells.Item(137, 10).Value = 40760
$ws.Cells.Item(137, 12).Value = 40760
$ws.Cells.Item(137, 14).Value = -50960

# --- BSM sheet: 26 cell update(s) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(8, 8).Value = 5636.7144
$ws.Cells.Item(8, 9).Value = 539.5
$ws.Cells.Item(8, 11).Value = 539.5
$ws.Cells.Item(8, 13).Value = -399.5
$ws.Cells.Item(59, 8).Value = 46890
$ws.Cells.Item(59, 10).Value = 46890
$ws.Cells.Item(59, 12).Value = 46890
$ws.Cells.Item(59, 14).Value = -48584
$ws.Cells.Item(107, 8).Value = 3515.8572
$ws.Cells.Item(107, 9).Value = 3227.75
$ws.Cells.Item(107, 10).Value = 3900
$ws.Cells.Item(107, 11).Value = 3227.75
$ws.Cells.Item(107, 12).Value = 3900
$ws.Cells.Item(107, 13).Value = -1307.75
$ws.Cells.Item(107, 14).Value = -7740
$ws.Cells.Item(134, 8).Value = 2064.5833
$ws.Cells.Item(134, 9).Value = 1502.174
$ws.Cells.Item(134, 10).Value = 15000
$ws.Cells.Item(134, 11).Value = 4506.522
$ws.Cells.Item(134, 12).Value = 45000
$ws.Cells.Item(134, 13).Value = -1971.522
$ws.Cells.Item(134, 14).Value = -50070
$ws.Cells.Item(137, 8).Value = 45300
$ws.Cells.Item(137, 10).Value = 45300
$ws.Cells.Item(137, 12).Value = 45300
$ws.Cells.Item(137, 14).Value = -55500

# --- CRP sheet: 52 cell update(s) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(50, 8).Value = 30124
$ws.Cells.Item(50, 10).Value = 30124
$ws.Cells.Item(50, 12).Value = 30124
$ws.Cells.Item(50, 14).Value = -31374
$ws.Cells.Item(51, 8).Value = 31514.4
$ws.Cells.Item(51, 10).Value = 31514.4
$ws.Cells.Item(51, 12).Value = 31514.4
$ws.Cells.Item(51, 14).Value = -32986.4
$ws.Cells.Item(58, 8).Value = 2950.672
$ws.Cells.Item(58, 9).Value = 1719.82
$ws.Cells.Item(58, 10).Value = 8545.454
$ws.Cells.Item(58, 11).Value = 1719.82
$ws.Cells.Item(58, 12).Value = 8545.454
$ws.Cells.Item(58, 13).Value = -1516.82
$ws.Cells.Item(58, 14).Value = -8951.454
$ws.Cells.Item(60, 8).Value = 18084.334
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = 18084.334
$ws.Cells.Item(60, 11).Value = 0
$ws.Cells.Item(60, 12).Value = 18084.334
$ws.Cells.Item(60, 13).Value = $null
$ws.Cells.Item(60, 14).Value = -19106.334
$ws.Cells.Item(61, 8).Value = 31514.4
$ws.Cells.Item(61, 10).Value = 31514.4
$ws.Cells.Item(61, 12).Value = 31514.4
$ws.Cells.Item(61, 14).Value = -32210.4
$ws.Cells.Item(87, 8).Value = 30466.666
$ws.Cells.Item(87, 10).Value = 30466.666
$ws.Cells.Item(87, 12).Value = 30466.666
$ws.Cells.Item(87, 14).Value = -32838.666
$ws.Cells.Item(90, 8).Value = 30466.666
$ws.Cells.Item(90, 10).Value = 30466.666
$ws.Cells.Item(90, 12).Value = 91399.99800000001
$ws.Cells.Item(90, 14).Value = -103255.998
$ws.Cells.Item(134, 8).Value = 8357.277
$ws.Cells.Item(134, 9).Value = 10177.583
$ws.Cells.Item(134, 10).Value = 4716.6665
$ws.Cells.Item(134, 11).Value = 30532.749
$ws.Cells.Item(134, 12).Value = 14149.9995
$ws.Cells.Item(134, 13).Value = -27997.749
$ws.Cells.Item(134, 14).Value = -19219.9995
$ws.Cells.Item(136, 8).Value = 2950.672
$ws.Cells.Item(136, 9).Value = 1719.82
$ws.Cells.Item(136, 10).Value = 8545.454
$ws.Cells.Item(136, 11).Value = 5159.46
$ws.Cells.Item(136, 12).Value = 25636.362
$ws.Cells.Item(136, 13).Value = -2609.46
$ws.Cells.Item(136, 14).Value = -30736.362
$ws.Cells.Item(137, 8).Value = 41818.57
$ws.Cells.Item(137, 10).Value = 41818.57
$ws.Cells.Item(137, 12).Value = 41818.57
$ws.Cells.Item(137, 14).Value = -52018.57

# --- CUL sheet: 18 cell update(s) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 711.25
$ws.Cells.Item(92, 9).Value = 712.8570999999999
$ws.Cells.Item(92, 10).Value = 700
$ws.Cells.Item(92, 11).Value = 2138.5713
$ws.Cells.Item(92, 12).Value = 2100
$ws.Cells.Item(92, 13).Value = -890.5712999999996
$ws.Cells.Item(92, 14).Value = -4596
$ws.Cells.Item(113, 8).Value = 3205709.8
$ws.Cells.Item(113, 9).Value = 587.14813
$ws.Cells.Item(113, 11).Value = 1761.44439
$ws.Cells.Item(113, 13).Value = 408.5556099999999
$ws.Cells.Item(131, 8).Value = 758.36365
$ws.Cells.Item(131, 9).Value = 253.7
$ws.Cells.Item(131, 10).Value = 815.06744
$ws.Cells.Item(131, 11).Value = 761.0999999999999
$ws.Cells.Item(131, 12).Value = 2445.20232
$ws.Cells.Item(131, 13).Value = 4278.9
$ws.Cells.Item(131, 14).Value = -12525.20232

# --- GSM sheet: 22 cell update(s) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(13, 8).Value = 1275.3334
$ws.Cells.Item(13, 9).Value = 854.1429000000001
$ws.Cells.Item(13, 10).Value = 2749.5
$ws.Cells.Item(13, 11).Value = 854.1429000000001
$ws.Cells.Item(13, 12).Value = 2749.5
$ws.Cells.Item(13, 13).Value = -715.1429000000001
$ws.Cells.Item(13, 14).Value = -3027.5
$ws.Cells.Item(46, 8).Value = 32417.334
$ws.Cells.Item(46, 10).Value = 34900.8
$ws.Cells.Item(46, 12).Value = 34900.8
$ws.Cells.Item(46, 14).Value = -35212.8
$ws.Cells.Item(132, 8).Value = 3301.6052
$ws.Cells.Item(132, 9).Value = 2398.0386
$ws.Cells.Item(132, 10).Value = 5259.3335
$ws.Cells.Item(132, 11).Value = 7194.1158
$ws.Cells.Item(132, 12).Value = 15778.0005
$ws.Cells.Item(132, 13).Value = -4664.1158
$ws.Cells.Item(132, 14).Value = -20838.0005
$ws.Cells.Item(137, 8).Value = 40454
$ws.Cells.Item(137, 10).Value = 40454
$ws.Cells.Item(137, 12).Value = 40454
$ws.Cells.Item(137, 14).Value = -50654

# --- LTW sheet: 25 cell update(s) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(117, 8).Value = 39892
$ws.Cells.Item(117, 10).Value = 39892
$ws.Cells.Item(117, 12).Value = 39892
$ws.Cells.Item(117, 14).Value = -49070
$ws.Cells.Item(122, 8).Value = 9088.200000000001
$ws.Cells.Item(122, 9).Value = 8659.333000000001
$ws.Cells.Item(122, 10).Value = 9272
$ws.Cells.Item(122, 11).Value = 25977.999
$ws.Cells.Item(122, 12).Value = 27816
$ws.Cells.Item(122, 13).Value = -23527.999
$ws.Cells.Item(122, 14).Value = -32716
$ws.Cells.Item(132, 8).Value = 8549.866
$ws.Cells.Item(132, 9).Value = 4095
$ws.Cells.Item(132, 10).Value = 17459.6
$ws.Cells.Item(132, 11).Value = 12285
$ws.Cells.Item(132, 12).Value = 52378.8
$ws.Cells.Item(132, 13).Value = -9755
$ws.Cells.Item(132, 14).Value = -57438.8
$ws.Cells.Item(136, 8).Value = 4482.6
$ws.Cells.Item(136, 9).Value = 1554.3334
$ws.Cells.Item(136, 10).Value = 8875
$ws.Cells.Item(136, 11).Value = 4663.0002
$ws.Cells.Item(136, 12).Value = 26625
$ws.Cells.Item(136, 13).Value = -2113.0002
$ws.Cells.Item(136, 14).Value = -31725

# --- WVR sheet: 21 cell update(s) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 7142.2104
$ws.Cells.Item(122, 9).Value = 5349.8
$ws.Cells.Item(122, 10).Value = 9133.777
$ws.Cells.Item(122, 11).Value = 16049.4
$ws.Cells.Item(122, 12).Value = 27401.331
$ws.Cells.Item(122, 13).Value = -13599.4
$ws.Cells.Item(122, 14).Value = -32301.331
$ws.Cells.Item(132, 8).Value = 7250296
$ws.Cells.Item(132, 9).Value = 3707.389
$ws.Cells.Item(132, 10).Value = 33338014
$ws.Cells.Item(132, 11).Value = 11122.167
$ws.Cells.Item(132, 12).Value = 100014042
$ws.Cells.Item(132, 13).Value = -8592.167000000001
$ws.Cells.Item(132, 14).Value = -100019102
$ws.Cells.Item(136, 8).Value = 13840.786
$ws.Cells.Item(136, 9).Value = 12310.333
$ws.Cells.Item(136, 10).Value = 16595.6
$ws.Cells.Item(136, 11).Value = 36930.999
$ws.Cells.Item(136, 12).Value = 49786.8
$ws.Cells.Item(136, 13).Value = -34380.999
$ws.Cells.Item(136, 14).Value = -54886.8

